$d = $word.ActiveDocument

# Locate the three target paragraphs by their (pre-edit) text content so the
# script is resilient to any paragraph renumbering.
$paraExploration = $null
$paraAcceptable = $null
$paraLast = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $t = $pp.Range.Text
    if ($t -like "Soar-RL*default exploration policy is epsilon-greedy*") {
        $paraExploration = $pp
    } elseif ($t -like "Acceptable values for epsilon are numbers between 0 and 1*") {
        $paraAcceptable = $pp
    } elseif ($t -like "With this explanation*this tutorial.*") {
        $paraLast = $pp
    }
}

# ------------------------------------------------------------------
# Paragraph 1: "Soar-RL's default exploration policy is epsilon-greedy, ..."
# ------------------------------------------------------------------
$r = $paraExploration.Range

$find = $r.Find
$find.ClearFormatting()
$find.Text = "Soar-RL's default exploration policy is "
$found = $find.Execute()
$r.Text = "When Soar is first started, the default exploration policy is "

$r.Collapse(0)
$r.InsertAfter("softmax")
$softmaxRange = $r.Duplicate

$r.Collapse(0)
$r.InsertAfter(".  However, the first time Soar-RL is enabled, the architecture automatically changes the exploration policy to ")

$softmaxRange.Font.Italic = 1

$r3 = $paraExploration.Range
$find2 = $r3.Find
$find2.ClearFormatting()
$find2.Text = "with an epsilon value of 0.1.  This policy states that 90% of the time the operator with greatest numerical preference value is chosen, while the remaining 10% of the time a random selection is made from all acceptable proposed operators.  You can change the epsilon value by issuing the following command:"
$found2 = $find2.Execute()

$r3.Text = "a policy more suitable for RL agents.  The default value of "

$r3.Collapse(0)
$r3.InsertAfter("epsilon")
$epsilon1 = $r3.Duplicate

$r3.Collapse(0)
$r3.InsertAfter(" is 0.1, dictating that 90% of the time the operator with greatest numerical preference value is chosen, while the remaining 10% of the time a random selection is made from all acceptable proposed operators.  You can change the ")

$r3.Collapse(0)
$r3.InsertAfter("epsilon")
$epsilon2 = $r3.Duplicate

$r3.Collapse(0)
$r3.InsertAfter(" value by issuing the following command:")

$epsilon1.Font.Italic = 1
$epsilon2.Font.Italic = 1

# ------------------------------------------------------------------
# Paragraph 2: "Acceptable values for epsilon are numbers between 0 and 1 ..."
# ------------------------------------------------------------------
$r4 = $paraAcceptable.Range
$find3 = $r4.Find
$find3.ClearFormatting()
$find3.Text = "Acceptable values for epsilon are numbers between 0 and 1"
$found3 = $find3.Execute()

$r4.Text = "Acceptable values for "
$r4.Collapse(0)
$r4.InsertAfter("epsilon")
$epsilon3 = $r4.Duplicate
$r4.Collapse(0)
$r4.InsertAfter(" are numbers between 0 and 1")

$epsilon3.Font.Italic = 1

# ------------------------------------------------------------------
# Add a trailing empty paragraph after the tutorial's closing paragraph.
# ------------------------------------------------------------------
$paraLast.Range.InsertParagraphAfter()
